$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 values
$ws.Range("O2").Value = 0.9267084217479559
$ws.Range("P2").Value = 0.9267084217479558
$ws.Range("S2").Value = 0.9267084217479559
$ws.Range("T2").Value = 0.9267084217479558

# Add new shared string value "Resolving-Mac" via row 3 column D
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Cort"
$ws.Range("C3").Value = "Sstr3"
$ws.Range("D3").Value = "Resolving-Mac"

$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1100146666666667
$ws.Range("H3").Value = 0.330044
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.02333733333333333
$ws.Range("N3").Value = 0.070012
$ws.Range("O3").Value = 0.07329157825204423
$ws.Range("P3").Value = 0.07329157825204423
$ws.Range("Q3").Value = 0.002567448947555556
$ws.Range("R3").Value = 0.023107040528
$ws.Range("S3").Value = 0.07329157825204423
$ws.Range("T3").Value = 0.07329157825204423
